$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11 for "Nason Creek Lower 08",
# shifting existing rows 11-16 down to 12-17.
$ws.Rows.Item(11).Insert()

# Populate the new row 11 (Nason Creek Lower 08)
$ws.Cells.Item(11, 1).Value = "Nason Creek Lower 08"
$ws.Cells.Item(11, 2).Value = "Wenatchee"
$ws.Cells.Item(11, 3).Value = "Lower Nason Creek"
$ws.Cells.Item(11, 4).Value = "yes"
$ws.Cells.Item(11, 5).Value = "yes"
$ws.Cells.Item(11, 6).Value = "yes"
$ws.Cells.Item(11, 7).Value = 3
$ws.Cells.Item(11, 9).Value = 3
$ws.Cells.Item(11, 10).Value = 5
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 3
$ws.Cells.Item(11, 13).Value = 3
$ws.Cells.Item(11, 14).Value = 1
$ws.Cells.Item(11, 15).Value = 1
$ws.Cells.Item(11, 16).Value = 3
$ws.Cells.Item(11, 17).Value = 1
$ws.Cells.Item(11, 18).Value = 2
$ws.Cells.Item(11, 19).Value = 1
$ws.Cells.Item(11, 20).Value = 22
$ws.Cells.Item(11, 21).Value = 0.4888888888888889
$ws.Cells.Item(11, 22).Value = 5
$ws.Cells.Item(11, 23).Value = 1
$ws.Cells.Item(11, 24).Value = "Off-Channel/Side-Channels,PoolQuantity&Quality,Riparian-Disturbance,Temperature-Rearing"
$ws.Cells.Item(11, 25).Value = "BankStability,Stability,Cover-Wood,Flow-SummerBaseFlow,FloodplainConnectivity,Riparian-CanopyCover,Riparian"
$ws.Cells.Item(11, 26).Value = "BankStability,Stability,Cover-Wood,Flow-SummerBaseFlow,FloodplainConnectivity,Off-Channel/Side-Channels,PoolQuantity&Quality,Riparian-Disturbance,Riparian-CanopyCover,Riparian,Temperature-Rearing"

# Update the unacceptable/at_risk summary columns (X, Y, Z) for all rows
# impacted by the new BankStability/ChannelStability/Riparian-CanopyCover/
# Riparian-Disturbance attribute columns.
$ws.Cells.Item(2, 24).Value = "BankStability,Cover-Wood,FloodplainConnectivity,Riparian-Disturbance"
$ws.Cells.Item(2, 25).Value = "ChannelStability,Stability,Flow-SummerBaseFlow,Off-Channel/Side-Channels,Riparian-CanopyCover,Riparian"
$ws.Cells.Item(2, 26).Value = "BankStability,ChannelStability,Stability,Cover-Wood,Flow-SummerBaseFlow,FloodplainConnectivity,Off-Channel/Side-Channels,Riparian-Disturbance,Riparian-CanopyCover,Riparian"
$ws.Cells.Item(3, 24).Value = "Cover-Wood,PoolQuantity&Quality"
$ws.Cells.Item(3, 25).Value = "Flow-SummerBaseFlow,Off-Channel/Side-Channels,Riparian-CanopyCover"
$ws.Cells.Item(3, 26).Value = "Cover-Wood,Flow-SummerBaseFlow,Off-Channel/Side-Channels,PoolQuantity&Quality,Riparian-CanopyCover"
$ws.Cells.Item(4, 25).Value = "BankStability,ChannelStability,Stability,Cover-Wood,Flow-SummerBaseFlow,FloodplainConnectivity,Off-Channel/Side-Channels,Riparian-Disturbance,Riparian-CanopyCover,Riparian,Temperature-Rearing"
$ws.Cells.Item(4, 26).Value = "BankStability,ChannelStability,Stability,Cover-Wood,Flow-SummerBaseFlow,FloodplainConnectivity,Off-Channel/Side-Channels,Riparian-Disturbance,Riparian-CanopyCover,Riparian,Temperature-Rearing"
$ws.Cells.Item(5, 25).Value = "BankStability,ChannelStability,Stability,Flow-SummerBaseFlow,FloodplainConnectivity,Off-Channel/Side-Channels,Riparian-Disturbance,Riparian-CanopyCover,Riparian,Temperature-Rearing"
$ws.Cells.Item(5, 26).Value = "BankStability,ChannelStability,Stability,Flow-SummerBaseFlow,FloodplainConnectivity,Off-Channel/Side-Channels,Riparian-Disturbance,Riparian-CanopyCover,Riparian,Temperature-Rearing"
$ws.Cells.Item(6, 24).Value = "Temperature-Rearing"
$ws.Cells.Item(6, 25).Value = "BankStability,Flow-SummerBaseFlow,FloodplainConnectivity,Riparian-CanopyCover"
$ws.Cells.Item(6, 26).Value = "BankStability,Flow-SummerBaseFlow,FloodplainConnectivity,Riparian-CanopyCover,Temperature-Rearing"
$ws.Cells.Item(7, 24).Value = "PoolQuantity&Quality,Temperature-Rearing"
$ws.Cells.Item(7, 25).Value = "BankStability,ChannelStability,Stability,Cover-Wood,FloodplainConnectivity,Off-Channel/Side-Channels,Riparian-Disturbance,Riparian-CanopyCover,Riparian"
$ws.Cells.Item(7, 26).Value = "BankStability,ChannelStability,Stability,Cover-Wood,FloodplainConnectivity,Off-Channel/Side-Channels,PoolQuantity&Quality,Riparian-Disturbance,Riparian-CanopyCover,Riparian,Temperature-Rearing"
$ws.Cells.Item(8, 24).Value = "PoolQuantity&Quality,Riparian-CanopyCover,Temperature-Rearing"
$ws.Cells.Item(8, 25).Value = "BankStability,ChannelStability,Stability,Cover-Wood,Flow-SummerBaseFlow,FloodplainConnectivity,Off-Channel/Side-Channels,Riparian-Disturbance,Riparian"
$ws.Cells.Item(8, 26).Value = "BankStability,ChannelStability,Stability,Cover-Wood,Flow-SummerBaseFlow,FloodplainConnectivity,Off-Channel/Side-Channels,PoolQuantity&Quality,Riparian-Disturbance,Riparian-CanopyCover,Riparian,Temperature-Rearing"
$ws.Cells.Item(9, 24).Value = "PoolQuantity&Quality,Temperature-Rearing"
$ws.Cells.Item(9, 25).Value = "BankStability,ChannelStability,Stability,Cover-Wood,Flow-SummerBaseFlow,FloodplainConnectivity,Off-Channel/Side-Channels,Riparian-Disturbance,Riparian-CanopyCover,Riparian"
$ws.Cells.Item(9, 26).Value = "BankStability,ChannelStability,Stability,Cover-Wood,Flow-SummerBaseFlow,FloodplainConnectivity,Off-Channel/Side-Channels,PoolQuantity&Quality,Riparian-Disturbance,Riparian-CanopyCover,Riparian,Temperature-Rearing"
$ws.Cells.Item(10, 24).Value = "PoolQuantity&Quality,Temperature-Rearing"
$ws.Cells.Item(10, 25).Value = "BankStability,ChannelStability,Stability,Cover-Wood,Flow-SummerBaseFlow,FloodplainConnectivity,Off-Channel/Side-Channels,Riparian-Disturbance,Riparian-CanopyCover,Riparian"
$ws.Cells.Item(10, 26).Value = "BankStability,ChannelStability,Stability,Cover-Wood,Flow-SummerBaseFlow,FloodplainConnectivity,Off-Channel/Side-Channels,PoolQuantity&Quality,Riparian-Disturbance,Riparian-CanopyCover,Riparian,Temperature-Rearing"
$ws.Cells.Item(12, 24).Value = "BankStability,Cover-Wood,FloodplainConnectivity,Off-Channel/Side-Channels,Riparian-CanopyCover,Temperature-Rearing"
$ws.Cells.Item(12, 25).Value = "ChannelStability,Stability,Flow-SummerBaseFlow,PoolQuantity&Quality,Riparian-Disturbance,Riparian"
$ws.Cells.Item(12, 26).Value = "BankStability,ChannelStability,Stability,Cover-Wood,Flow-SummerBaseFlow,FloodplainConnectivity,Off-Channel/Side-Channels,PoolQuantity&Quality,Riparian-Disturbance,Riparian-CanopyCover,Riparian,Temperature-Rearing"
$ws.Cells.Item(13, 24).Value = "BankStability,Cover-Wood,FloodplainConnectivity,Off-Channel/Side-Channels,Riparian-CanopyCover,Temperature-Rearing"
$ws.Cells.Item(13, 25).Value = "ChannelStability,Stability,Flow-SummerBaseFlow,PoolQuantity&Quality,Riparian-Disturbance,Riparian"
$ws.Cells.Item(13, 26).Value = "BankStability,ChannelStability,Stability,Cover-Wood,Flow-SummerBaseFlow,FloodplainConnectivity,Off-Channel/Side-Channels,PoolQuantity&Quality,Riparian-Disturbance,Riparian-CanopyCover,Riparian,Temperature-Rearing"
$ws.Cells.Item(14, 24).Value = "PoolQuantity&Quality,Riparian-Disturbance,Temperature-Rearing"
$ws.Cells.Item(14, 25).Value = "BankStability,ChannelStability,Stability,CoarseSubstrate,Cover-Wood,Flow-SummerBaseFlow,FloodplainConnectivity,Off-Channel/Side-Channels,Riparian-CanopyCover,Riparian"
$ws.Cells.Item(14, 26).Value = "BankStability,ChannelStability,Stability,CoarseSubstrate,Cover-Wood,Flow-SummerBaseFlow,FloodplainConnectivity,Off-Channel/Side-Channels,PoolQuantity&Quality,Riparian-Disturbance,Riparian-CanopyCover,Riparian,Temperature-Rearing"
$ws.Cells.Item(15, 24).Value = "Cover-Wood,Riparian-Disturbance,Riparian-CanopyCover,Riparian,Temperature-Rearing"
$ws.Cells.Item(15, 25).Value = "BankStability,ChannelStability,Stability,CoarseSubstrate,Flow-SummerBaseFlow,FloodplainConnectivity,Off-Channel/Side-Channels,PoolQuantity&Quality"
$ws.Cells.Item(15, 26).Value = "BankStability,ChannelStability,Stability,CoarseSubstrate,Cover-Wood,Flow-SummerBaseFlow,FloodplainConnectivity,Off-Channel/Side-Channels,PoolQuantity&Quality,Riparian-Disturbance,Riparian-CanopyCover,Riparian,Temperature-Rearing"
$ws.Cells.Item(16, 24).Value = "Cover-Wood,Riparian-Disturbance,Riparian-CanopyCover,Riparian,Temperature-Rearing"
$ws.Cells.Item(16, 25).Value = "BankStability,ChannelStability,Stability,CoarseSubstrate,Flow-SummerBaseFlow,FloodplainConnectivity,Off-Channel/Side-Channels,PoolQuantity&Quality"
$ws.Cells.Item(16, 26).Value = "BankStability,ChannelStability,Stability,CoarseSubstrate,Cover-Wood,Flow-SummerBaseFlow,FloodplainConnectivity,Off-Channel/Side-Channels,PoolQuantity&Quality,Riparian-Disturbance,Riparian-CanopyCover,Riparian,Temperature-Rearing"
$ws.Cells.Item(17, 24).Value = "Temperature-Rearing"
$ws.Cells.Item(17, 25).Value = "BankStability,ChannelStability,Stability,CoarseSubstrate,Cover-Wood,FloodplainConnectivity,Off-Channel/Side-Channels,PoolQuantity&Quality,Riparian-Disturbance,Riparian"
$ws.Cells.Item(17, 26).Value = "BankStability,ChannelStability,Stability,CoarseSubstrate,Cover-Wood,FloodplainConnectivity,Off-Channel/Side-Channels,PoolQuantity&Quality,Riparian-Disturbance,Riparian,Temperature-Rearing"
